$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The player name cell in the original workbook ends with a non-breaking
# space (U+00A0) rather than a plain space - reuse the exact text already
# present in F2 so every row matches byte-for-byte.
$playerName = $ws.Cells.Item(2, 6).Value2

$data = @(
    @(" Nov 2 2020", " Abu Dhabi", "Capitals won by 6 wickets (with 6 balls remaining)", "Delhi Capitals", "Royal Challengers Bangalore", $playerName, "60", "46", "5", "1", "130.43"),
    @(" Nov 10 2020", " Dubai (DSC)", "Mumbai won by 5 wickets (with 8 balls remaining)", "Delhi Capitals", "Mumbai Indians", $playerName, "2", "4", "0", "0", "50.00"),
    @(" Oct 17 2020", " Sharjah", "Capitals won by 5 wickets (with 1 ball remaining)", "Delhi Capitals", "Chennai Super Kings", $playerName, "8", "10", "1", "0", "80.00"),
    @(" Oct 14 2020", " Dubai (DSC)", "Capitals won by 13 runs", "Delhi Capitals", "Rajasthan Royals", $playerName, "2", "9", "0", "0", "22.22"),
    @(" Oct 24 2020", " Abu Dhabi", "KKR won by 59 runs", "Delhi Capitals", "Kolkata Knight Riders", $playerName, "0", "1", "0", "0", "0.00"),
    @(" Oct 11 2020", " Abu Dhabi", "Mumbai won by 5 wickets (with 2 balls remaining)", "Delhi Capitals", "Mumbai Indians", $playerName, "15", "15", "3", "0", "100.00"),
    @(" Nov 5 2020", " Dubai (DSC)", "Mumbai won by 57 runs", "Delhi Capitals", "Mumbai Indians", $playerName, "0", "3", "0", "0", "0.00"),
    @(" Oct 27 2020", " Dubai (DSC)", "Sunrisers won by 88 runs", "Delhi Capitals", "Sunrisers Hyderabad", $playerName, "26", "19", "3", "1", "136.84")
)

# Ensure the whole data range keeps values as text, matching the
# "numberStoredAsText" semantics of the original workbook.
$dataRange = $ws.Range("A2:K9")
$dataRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($rowIndex, $j + 1).Value = $rowData[$j]
    }
}
